$d = $word.ActiveDocument

# Locate the paragraph that ends with "...Resolvable URNs." and add the
# new list item right after it.
$r = $d.Content
$found = $r.Find.Execute("Resources Embeddings Registry / Index, Naming (prompts, placeholders). Resolvable URNs.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorPara = $r.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$newPara.Range.Text = "URN ID: Contexts. FCA Contexts Prime ID Embedding."
